$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The A8 label loses its leading space: " 900 - 1000" -> "900 - 1000"
$ws.Range("A8").Value = "900 - 1000"

# Update the active selection to A8 (was B8)
$ws.Range("A8").Select()
